# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "67.863.86" },
    @{ Cell = "E2"; Value = "  +1.14%  " },
    @{ Cell = "D3"; Value = "3.515.18" },
    @{ Cell = "E3"; Value = "  +0.16%  " },
    @{ Cell = "E4"; Value = "  -0.01%  " },
    @{ Cell = "D5"; Value = "600.43" },
    @{ Cell = "E5"; Value = "  +0.91%  " },
    @{ Cell = "D6"; Value = "181.12" },
    @{ Cell = "E6"; Value = "  +4.43%  " },
    @{ Cell = "E7"; Value = "  +0.02%  " },
    @{ Cell = "D8"; Value = "3.515.22" },
    @{ Cell = "E8"; Value = "  +0.17%  " },
    @{ Cell = "D9"; Value = "0.595" },
    @{ Cell = "E9"; Value = "  -0.15%  " },
    @{ Cell = "D10"; Value = "0.140" },
    @{ Cell = "E10"; Value = "  +6.60%  " },
    @{ Cell = "D11"; Value = "7.15" },
    @{ Cell = "E11"; Value = "  -1.93%  " },
    @{ Cell = "E12"; Value = "  +0.65%  " },
    @{ Cell = "D13"; Value = "4.119.61" },
    @{ Cell = "E13"; Value = "  +0.08%  " },
    @{ Cell = "D14"; Value = "32.68" },
    @{ Cell = "E14"; Value = "  +12.27%  " },
    @{ Cell = "E15"; Value = "  +1.14%  " },
    @{ Cell = "D16"; Value = "67.856.00" },
    @{ Cell = "E16"; Value = "  +1.11%  " },
    @{ Cell = "E17"; Value = "  +0.18%  " },
    @{ Cell = "D18"; Value = "3.512.31" },
    @{ Cell = "E18"; Value = "  -0.32%  " },
    @{ Cell = "D19"; Value = "6.35" },
    @{ Cell = "E19"; Value = "  +0.38%  " },
    @{ Cell = "D20"; Value = "14.47" },
    @{ Cell = "D21"; Value = "399.79" },
    @{ Cell = "E21"; Value = "  +1.00%  " },
    @{ Cell = "D22"; Value = "7.98" },
    @{ Cell = "E22"; Value = "  -0.86%  " },
    @{ Cell = "D23"; Value = "73.67" },
    @{ Cell = "E23"; Value = "  +0.70%  " },
    @{ Cell = "E24"; Value = "  +0.93%  " },
    @{ Cell = "D25"; Value = "0.999" },
    @{ Cell = "E25"; Value = "  -0.08%  " },
    @{ Cell = "E26"; Value = "  +0.37%  " },
    @{ Cell = "E27"; Value = "  +1.01%  " },
    @{ Cell = "D28"; Value = "10.51" },
    @{ Cell = "E28"; Value = "  +2.01%  " },
    @{ Cell = "E29"; Value = "  -2.60%  " },
    @{ Cell = "E30"; Value = "  +0.01%  " },
    @{ Cell = "D31"; Value = "6.24" },
    @{ Cell = "E31"; Value = "  -1.29%  " },
    @{ Cell = "E32"; Value = "  -0.25%  " },
    @{ Cell = "E33"; Value = "  +1.46%  " },
    @{ Cell = "D34"; Value = "23.89" },
    @{ Cell = "E34"; Value = "  +0.21%  " },
    @{ Cell = "D35"; Value = "7.49" },
    @{ Cell = "E35"; Value = "  +1.47%  " },
    @{ Cell = "E36"; Value = "  +0.28%  " },
    @{ Cell = "D37"; Value = "1.62" },
    @{ Cell = "E37"; Value = "  -3.74%  " },
    @{ Cell = "D38"; Value = "163.23" },
    @{ Cell = "E38"; Value = "  -0.03%  " },
    @{ Cell = "D39"; Value = "0.881" },
    @{ Cell = "E39"; Value = "  -0.46%  " },
    @{ Cell = "E40"; Value = "  +0.57%  " },
    @{ Cell = "B41"; Value = "RenderToken" },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" },
    @{ Cell = "D41"; Value = "6.99" },
    @{ Cell = "E41"; Value = "  -1.18%  " },
    @{ Cell = "B42"; Value = "dogwifhat" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" },
    @{ Cell = "D42"; Value = "2.76" },
    @{ Cell = "E42"; Value = "  +6.69%  " },
    @{ Cell = "D43"; Value = "2.889.46" },
    @{ Cell = "E43"; Value = "  +2.59%  " },
    @{ Cell = "D44"; Value = "4.69" },
    @{ Cell = "E44"; Value = "  -0.22%  " },
    @{ Cell = "D45"; Value = "0.0735" },
    @{ Cell = "E45"; Value = "  -1.75%  " },
    @{ Cell = "D46"; Value = "26.26" },
    @{ Cell = "E46"; Value = "  -0.67%  " },
    @{ Cell = "D47"; Value = "26.89" },
    @{ Cell = "E47"; Value = "  -2.41%  " },
    @{ Cell = "D48"; Value = "42.38" },
    @{ Cell = "E48"; Value = "  -0.79%  " },
    @{ Cell = "B49"; Value = "VeChain" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" },
    @{ Cell = "D49"; Value = "0.0303" },
    @{ Cell = "E49"; Value = "  -0.38%  " },
    @{ Cell = "B50"; Value = "Bittensor" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" },
    @{ Cell = "D50"; Value = "344.78" },
    @{ Cell = "E50"; Value = "  +1.34%  " },
    @{ Cell = "E51"; Value = "  -1.56%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
